# Weekly update: insert the two most-recent price records for this
# commodity/market subset. This pushes the existing rows 135..214 down to
# 137..216 and fills the two newly-opened rows (135, 136) with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 135:136 (shifts old rows 135-214 down to 137-216).
$ws.Range("A135:A136").EntireRow.Insert()

# --- New row 135 -----------------------------------------------------------
$ws.Cells.Item(135, 1).Value  = 1
$ws.Cells.Item(135, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(135, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(135, 4).Value  = 45089
$ws.Cells.Item(135, 5).Value  = 15
$ws.Cells.Item(135, 6).Value  = 100114001
$ws.Cells.Item(135, 7).Value  = "Papa"
$ws.Cells.Item(135, 8).Value  = "Asterix"
$ws.Cells.Item(135, 9).Value  = "1a (guarda)"
$ws.Cells.Item(135, 10).Value = 1000
$ws.Cells.Item(135, 11).Value = 12000
$ws.Cells.Item(135, 12).Value = 13000
$ws.Cells.Item(135, 13).Value = 12500
$ws.Cells.Item(135, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(135, 15).Value = "Región Metropolitana"
$ws.Cells.Item(135, 16).Value = 500
$ws.Cells.Item(135, 17).Value = 25
$ws.Cells.Item(135, 18).Value = "Hortaliza"

# --- New row 136 -----------------------------------------------------------
$ws.Cells.Item(136, 1).Value  = 1
$ws.Cells.Item(136, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(136, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(136, 4).Value  = 45089
$ws.Cells.Item(136, 5).Value  = 15
$ws.Cells.Item(136, 6).Value  = 100114001
$ws.Cells.Item(136, 7).Value  = "Papa"
$ws.Cells.Item(136, 8).Value  = "Cardinal"
$ws.Cells.Item(136, 9).Value  = "1a nueva(o)"
$ws.Cells.Item(136, 10).Value = 1000
$ws.Cells.Item(136, 11).Value = 14000
$ws.Cells.Item(136, 12).Value = 15000
$ws.Cells.Item(136, 13).Value = 14500
$ws.Cells.Item(136, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(136, 15).Value = "Región Metropolitana"
$ws.Cells.Item(136, 16).Value = 580
$ws.Cells.Item(136, 17).Value = 25
$ws.Cells.Item(136, 18).Value = "Hortaliza"
